$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.208.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6989"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07730"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3046"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7174"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.149"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.225.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.753"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007718"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9990"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.109.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1481"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.006"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.048"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.418"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.434"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.482"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05192"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.165"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7082"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9993"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.723"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9374"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.140.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4280"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.892"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9996"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.792"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.008.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.159"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.98%  "

